# Add a new row (82) of data to the first worksheet, matching the
# pattern of the existing date/value rows (A = serial date, B = value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 82

# Copy the style of the previous date cell (A81) onto the new cell (A82)
# so that the date number format (s="2") carries over.
$ws.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = 45884
$ws.Cells.Item($newRow, 2).Value = 0.06654624964350926
